# final draft text & metadata edits before review
#
# Updates the Personnel sheet: corrects a contributor's given/sur name,
# and leaves the sheet scrolled/selected the way the author left it
# (selection moved from G7 to C8, view scrolled down a few rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Personnel")

# --- Text corrections (row 7: givenName / surName) ---
$ws.Range("A7").Value = "S. Alejandra"
$ws.Range("C7").Value = "Casillo Cieza"

# --- View / selection state ---
$ws.Activate()

# Scroll the window so row 4 is at the top-left of the viewport, then move
# the selection to C8 (matching the author's final cursor position).
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C8").Select()
